# Apply the automatic-update edit: rows 3,4,5,6,7,9 in the active sheet
# had their species/record data cyclically reassigned among each other
# (A, B, D, E, F, G, H, Q, R columns, plus Y/AA start & end dates for
# rows 4 and 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("A3").Value = 112181500
$ws.Range("B3").Value = 89517
$ws.Range("E3").Value = 5447
$ws.Range("F3").Value = "Vedticka"
$ws.Range("G3").Value = "Fuscoporia viticola"
$ws.Range("H3").Value = "(Schwein.) Murrill"
$ws.Range("Q3").Value = 772346
$ws.Range("R3").Value = 7120286

# --- Row 4 ---
$ws.Range("A4").Value = 112181511
$ws.Range("Q4").Value = 772359
$ws.Range("R4").Value = 7120174
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-09-07"
$ws.Range("Y4").Style = "Normal"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-09-07"
$ws.Range("AA4").Style = "Normal"

# --- Row 5 ---
$ws.Range("A5").Value = 112181514
$ws.Range("Q5").Value = 772353
$ws.Range("R5").Value = 7120281

# --- Row 6 ---
$ws.Range("A6").Value = 112181532
$ws.Range("B6").Value = 89499
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 112
$ws.Range("F6").Value = "Stjärntagging"
$ws.Range("G6").Value = "Asterodon ferruginosus"
$ws.Range("H6").Value = "Pat."
$ws.Range("Q6").Value = 772340
$ws.Range("R6").Value = 7120223

# --- Row 7 ---
$ws.Range("A7").Value = 112181512
$ws.Range("B7").Value = 89553
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = "Ullticka"
$ws.Range("G7").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H7").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q7").Value = 772413
$ws.Range("R7").Value = 7120316

# --- Row 9 ---
$ws.Range("A9").Value = 112182926
$ws.Range("B9").Value = 5113
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 100526
$ws.Range("F9").Value = "Bronshjon"
$ws.Range("G9").Value = "Callidium coriaceum"
$ws.Range("H9").Value = "Paykull, 1800"
$ws.Range("Q9").Value = 772357
$ws.Range("R9").Value = 7120234
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2023-09-06"
$ws.Range("Y9").Style = "Normal"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2023-09-06"
$ws.Range("AA9").Style = "Normal"
